# feat: add 2022-Q1 data
#
# - Inserts a new "2022-Q1" worksheet (fund holdings for the quarter)
#   positioned between "2021-Q4" and "总计".
# - Prepends a matching summary row to the "总计" (totals) sheet.

function Set-TextValue {
    # Writes $value to $cell as TEXT, even when it looks like a number
    # (e.g. "010923", "2.30"), without leaving a stray number-format
    # style behind on the cell.
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right after "2021-Q4" (before "总计")
# ------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q1"

# Copy the header row formatting (style s=2) from the previous quarter
# sheet, then overwrite the header text for this sheet.
$afterSheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Copy the index-column formatting (style s=2) for A2, then fill in the
# single fund entry for 2022-Q1.
$afterSheet.Cells.Item(2, 1).Copy($newSheet.Cells.Item(2, 1))
$newSheet.Cells.Item(2, 1).Value = 0

Set-TextValue $newSheet.Cells.Item(2, 2) "010923"
$newSheet.Cells.Item(2, 3).Value = "永赢鑫欣混合"
Set-TextValue $newSheet.Cells.Item(2, 4) "2.30"
Set-TextValue $newSheet.Cells.Item(2, 5) "30.35"
Set-TextValue $newSheet.Cells.Item(2, 6) "3.30"
Set-TextValue $newSheet.Cells.Item(2, 7) "0.0759"
$newSheet.Cells.Item(2, 8).Value = 2

# ------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new top data row for
#    2022-Q1 and shift the existing quarters down by one row.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The Insert() above copies row-1 (header) formatting onto the new
# row 2 - clear that back to the unstyled data-row look first.
$totalSheet.Range("A2:D2").Style = "Normal"

# Restore the index-column style (s=2) on A2 by copying it from A3.
$totalSheet.Cells.Item(3, 1).Copy($totalSheet.Cells.Item(2, 1))

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 1
$totalSheet.Cells.Item(2, 4).Value = 0.08

# The row-insert also shifted the leading index column (A) down with
# everything else; the source data keeps each row's index equal to
# (row number - 2), so re-pin A3:A7 back to that sequence.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(7, 1).Value = 5
